# Test data update: rename the "Loading Details Name" value used for the
# 24V PSU load rows from "Main Processor 24V (A)" to "24V Rail(A)".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Column G (rows 8-13) holds the "Loading Details Name" shared-string value.
$ws.Range("G8:G13").Value = "24V Rail(A)"

# Reflect the new selection left behind after the edit.
$ws.Range("G8").Select()
